# Bump the "Förändrad" (Changed) date in column C for rows 2-66
# from serial 45171 (2023-09-02) to serial 45172 (2023-09-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C66").Value = 45172
